$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (deg columns) - B/D now mirror the old "16" column,
# C/E now mirror the old "20" column (Lichtwark passive tweak / deleted cols)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON)
$ws.Range("B2").Value = 432.99391268812496
$ws.Range("C2").Value = 514.43989045750004
$ws.Range("D2").Value = 432.99391268812496
$ws.Range("E2").Value = 514.43989045750004

# Row 3 (STR)
$ws.Range("B3").Value = 432.99391268812496
$ws.Range("C3").Value = 516.97407711000005
$ws.Range("D3").Value = 432.99391268812496
$ws.Range("E3").Value = 516.97407711000005

# Selection now only spans the edited block instead of the whole table
[void]$ws.Range("B1:E3").Select()
